# tabla dns completa (con valores no reales)
#
# The "Direcciones necesarias" column (E) is no longer needed on screen,
# so it gets selected and hidden - mirroring what Excel records when a
# user right-clicks the column header and chooses "Hide":
#   - column E's <col .../> gains hidden="1"
#   - the sheet's active selection becomes the full column (E1:E1048576)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns("E")
$col.Select()
$col.Hidden = $true
